$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All simulations")

# C7: "None." -> "None"
$ws.Range("C7").Value2 = "None"

# G9/G10: "3 from above " -> existing migration-rates string (shared with G7)
$migRates = "0; 0.001; 0.0025; 0.005; 0.0075; 0.01; 0.025; 0.05; 0.1; 0.2"
$ws.Range("G9").Value2 = $migRates
$ws.Range("G10").Value2 = $migRates

# I9/I10: simulation counts updated
$ws.Range("I9").Value2 = 30
$ws.Range("I10").Value2 = 90

# J9: "Not done" -> "Running"
$ws.Range("J9").Value2 = "Running"

# K8: new cell with "Not interesting"
$ws.Range("K8").Value2 = "Not interesting"

# Rows 9 & 10 grow taller to fit the updated migration-rates text
$ws.Rows.Item(9).RowHeight = 48
$ws.Rows.Item(10).RowHeight = 48

# Update selection to reflect J10 being the active cell
$ws.Range("J10").Select()
